$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new season row (row 27), reusing the formatting of the row above
$ws.Range("A26:E26").Copy()
$ws.Range("A27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Write D27 first as quoted text so it is stored as a shared string
# (its column is date-formatted, so a plain numeric value would be
# interpreted as a number) - then reapply D26's exact style.
$ws.Range("D27").Value = "'10.2"
$ws.Range("D26").Copy()
$ws.Range("D27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("E27").Value = "Start Dandelion Journey"
$ws.Range("A27").Value = "M4_02 Love 2022"
$ws.Range("B27").Value = 44600
$ws.Range("C27").Value = 44628

$ws.Range("A27").Select()
